$d = $word.ActiveDocument

# The paragraph carries a zero-width "_GoBack" bookmark sitting right at
# the end of the run text (a leftover of the last edit position). If we
# leave it in place, InsertXML's range-replace logic anchors the
# bookmark's start marker to the left edge of whatever new content lands
# at that same offset, which would yank <w:bookmarkStart/> in front of
# our new runs. Remove it first and re-create it (in the right spot, via
# the replacement markup below) so it ends up exactly where it started:
# right after the final run, before the paragraph mark.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# Find the exact run of text that needs to be split into
# spell/grammar-checked runs with interspersed <w:proofErr/> markers
# (mirroring what Word's background proofer emits while typing) and
# extended with the new trailing sentence.
$rng = $d.Content.Duplicate
$null = $rng.Find.Execute("Dsfsdaf fasds dfsdf  dasdf")

$inner = '<w:proofErr w:type="spellStart"/>' +
         '<w:r><w:t>Dsfsdaf</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/>' +
         '<w:r><w:t>fasds</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/>' +
         '<w:proofErr w:type="gramStart"/>' +
         '<w:r><w:t>dfsdf</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r><w:t xml:space="preserve">  </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/>' +
         '<w:r><w:t>dasdf</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:proofErr w:type="gramEnd"/>' +
         '<w:r><w:t>. This is new</w:t></w:r>' +
         '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
         '<w:bookmarkEnd w:id="0"/>'

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
       '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body><w:p>' + $inner + '</w:p></w:body>' +
       '</w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'

$rng.InsertXML($xml)
